$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 2038.26
$ws.Range("F13").Value = 41
$ws.Range("G13").Value = 8935.129999999999
$ws.Range("B14").Value = 17013.91
$ws.Range("F56").Value = 29
$ws.Range("G56").Value = 528.96
$ws.Range("F57").Value = 148
$ws.Range("G57").Value = 13843.92
$ws.Range("B71").Value = 77268.17
$ws.Range("F108").Value = 48
$ws.Range("G108").Value = 4570.08
$ws.Range("F117").Value = 333
$ws.Range("G117").Value = 21212.1
$ws.Range("F122").Value = 348
$ws.Range("G122").Value = 5755.92
$ws.Range("F129").Value = 205
$ws.Range("G129").Value = 21004.3
$ws.Range("F140").Value = 49
$ws.Range("G140").Value = 2190.3
$ws.Range("B143").Value = 347159.03
$ws.Range("F212").Value = 22
$ws.Range("G212").Value = 2533.96
$ws.Range("B216").Value = 6850.54
$ws.Range("F258").Value = 72
$ws.Range("G258").Value = 1825.92
$ws.Range("F260").Value = 10
$ws.Range("G260").Value = 782.6
$ws.Range("B264").Value = 19166.07
$ws.Range("F299").Value = 32
$ws.Range("G299").Value = 2471.04
$ws.Range("B303").Value = 23751.96
$ws.Range("F343").Value = 6
$ws.Range("G343").Value = 839.7
$ws.Range("B386").Value = 152518.76
$ws.Range("F390").Value = 94
$ws.Range("G390").Value = 20324.68
$ws.Range("B391").Value = 57077
$ws.Range("D391").Value = 93.08
$ws.Range("E391").Value = 111.2
$ws.Range("F391").Value = 1
$ws.Range("G391").Value = 93.08
$ws.Range("B392").Value = 61610
$ws.Range("D392").Value = 102.71
$ws.Range("E392").Value = 122.71
$ws.Range("F392").Value = 383
$ws.Range("G392").Value = 39337.93
$ws.Range("F409").Value = 183
$ws.Range("G409").Value = 25083.81
$ws.Range("F412").Value = 1
$ws.Range("G412").Value = 86.7
$ws.Range("F425").Value = 269
$ws.Range("G425").Value = 30117.24
$ws.Range("F434").Value = 59
$ws.Range("G434").Value = 6141.31
$ws.Range("F438").Value = 263
$ws.Range("G438").Value = 26589.3
$ws.Range("F444").Value = 66
$ws.Range("G444").Value = 6858.06
$ws.Range("F451").Value = 697
$ws.Range("G451").Value = 15117.93
$ws.Range("F452").Value = 263
$ws.Range("G452").Value = 1580.63
$ws.Range("F453").Value = 20
$ws.Range("G453").Value = 1662
$ws.Range("F457").Value = 47
$ws.Range("G457").Value = 8741.530000000001
$ws.Range("F467").Value = 25
$ws.Range("G467").Value = 508
$ws.Range("B471").Value = 635275.35
$ws.Range("F473").Value = 31
$ws.Range("G473").Value = 4755.71
$ws.Range("B486").Value = 60360.42
$ws.Range("F540").Value = 236
$ws.Range("G540").Value = 22797.6
$ws.Range("B546").Value = 74858.50999999999
$ws.Range("F555").Value = 11
$ws.Range("G555").Value = 1992.76
$ws.Range("F563").Value = 118
$ws.Range("G563").Value = 5593.2
$ws.Range("B575").Value = 79839.86
$ws.Range("F593").Value = 699
$ws.Range("G593").Value = 8954.190000000001
$ws.Range("F595").Value = 301
$ws.Range("G595").Value = 5938.73
$ws.Range("F596").Value = 391
$ws.Range("G596").Value = 6424.13
$ws.Range("F601").Value = 193
$ws.Range("G601").Value = 3755.78
$ws.Range("F607").Value = 338
$ws.Range("G607").Value = 5553.34
$ws.Range("B609").Value = 123356.72
$ws.Range("F674").Value = 634
$ws.Range("G674").Value = 4342.9
$ws.Range("F677").Value = 562
$ws.Range("G677").Value = 3765.4
$ws.Range("B681").Value = 49953.22
$ws.Range("F685").Value = 23
$ws.Range("G685").Value = 4344.93
$ws.Range("B689").Value = 6780.05
$ws.Range("F714").Value = 113
$ws.Range("G714").Value = 11130.5
$ws.Range("F720").Value = 179
$ws.Range("G720").Value = 5876.57
$ws.Range("B722").Value = 77646.21000000001
$ws.Range("F728").Value = 121
$ws.Range("G728").Value = 2675.31
$ws.Range("B743").Value = 14379.89
$ws.Range("F746").Value = 39
$ws.Range("G746").Value = 11148.93
$ws.Range("F759").Value = 13
$ws.Range("G759").Value = 1733.29
$ws.Range("F761").Value = 84
$ws.Range("G761").Value = 6178.2
$ws.Range("B765").Value = 91780.14
$ws.Range("F772").Value = 173
$ws.Range("G772").Value = 4705.6
$ws.Range("B774").Value = 82889.75999999999
$ws.Range("F800").Value = 465
$ws.Range("G800").Value = 7374.9
$ws.Range("F801").Value = 137
$ws.Range("G801").Value = 4536.07
$ws.Range("F803").Value = 117
$ws.Range("G803").Value = 3873.87
$ws.Range("F807").Value = 89
$ws.Range("G807").Value = 3894.64
$ws.Range("B808").Value = 54287.13
$ws.Range("F816").Value = 31
$ws.Range("G816").Value = 1918.9
$ws.Range("B830").Value = 33290.74
$ws.Range("F837").Value = 53
$ws.Range("G837").Value = 4862.75
$ws.Range("F838").Value = 46
$ws.Range("G838").Value = 3933
$ws.Range("B843").Value = 28849.35
$ws.Range("F882").Value = 13
$ws.Range("G882").Value = 8102.64
$ws.Range("B884").Value = 62764.23
$ws.Range("F889").Value = 18
$ws.Range("G889").Value = 1991.34
$ws.Range("F893").Value = 2
$ws.Range("G893").Value = 74.56
$ws.Range("F902").Value = 138
$ws.Range("G902").Value = 19872
$ws.Range("F903").Value = 265
$ws.Range("G903").Value = 31988.15
$ws.Range("B905").Value = 128284.02
$ws.Range("F922").Value = 20
$ws.Range("G922").Value = 2855.4
$ws.Range("F927").Value = 191
$ws.Range("G927").Value = 19655.81
$ws.Range("F929").Value = 35
$ws.Range("G929").Value = 2428.3
$ws.Range("F932").Value = 21
$ws.Range("G932").Value = 668.01
$ws.Range("F936").Value = 19
$ws.Range("G936").Value = 1082.81
$ws.Range("B937").Value = 75137.35000000001
$ws.Range("F940").Value = 134
$ws.Range("G940").Value = 5011.6
$ws.Range("F941").Value = 61
$ws.Range("G941").Value = 1171.81
$ws.Range("F942").Value = 87
$ws.Range("G942").Value = 569.85
$ws.Range("F943").Value = 245
$ws.Range("G943").Value = 9163
$ws.Range("F945").Value = 227
$ws.Range("G945").Value = 8489.799999999999
$ws.Range("B946").Value = 31595.3
$ws.Range("F992").Value = 5
$ws.Range("G992").Value = 825.15
$ws.Range("B1001").Value = 10923.12
$ws.Range("F1005").Value = 296
$ws.Range("G1005").Value = 22833.44
$ws.Range("B1009").Value = 562007.02
$ws.Range("B1016").Value = 4491955.79
$ws.Range("B1017").Value = 4491955.79
